$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price (D) and Volume (E) columns so that
# numeric-looking strings (e.g. "172.21", "11.00", "0.0372") are stored
# verbatim as text instead of being coerced into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Row 36 / 37 swap: Dai and Hedera exchange rank positions with updated values ---
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.101"
$ws.Range("E36").Value = "  -11.31%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.17%  "

# --- Remaining per-row price / volume updates ---
$ws.Range("D2").Value = "62.683.43"
$ws.Range("E2").Value = "  -8.43%  "
$ws.Range("D3").Value = "3.221.20"
$ws.Range("E3").Value = "  -10.77%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "172.21"
$ws.Range("E5").Value = "  -16.64%  "
$ws.Range("D6").Value = "503.48"
$ws.Range("E6").Value = "  -11.94%  "
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  -4.92%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "3.217.58"
$ws.Range("E9").Value = "  -10.69%  "
$ws.Range("D10").Value = "0.602"
$ws.Range("E10").Value = "  -12.46%  "
$ws.Range("D11").Value = "55.64"
$ws.Range("E11").Value = "  -13.33%  "
$ws.Range("E12").Value = "  -14.94%  "
$ws.Range("E13").Value = "  -12.58%  "
$ws.Range("D14").Value = "8.93"
$ws.Range("E14").Value = "  -13.88%  "
$ws.Range("D15").Value = "3.725.15"
$ws.Range("E15").Value = "  -10.75%  "
$ws.Range("E16").Value = "  -7.13%  "
$ws.Range("D17").Value = "3.210.88"
$ws.Range("E17").Value = "  -10.99%  "
$ws.Range("D18").Value = "62.479.21"
$ws.Range("E18").Value = "  -8.45%  "
$ws.Range("D19").Value = "16.99"
$ws.Range("D20").Value = "10.67"
$ws.Range("E20").Value = "  -13.24%  "
$ws.Range("D21").Value = "0.923"
$ws.Range("E21").Value = "  -13.88%  "
$ws.Range("D22").Value = "362.44"
$ws.Range("E22").Value = "  -10.97%  "
$ws.Range("D23").Value = "78.32"
$ws.Range("E23").Value = "  -7.93%  "
$ws.Range("D24").Value = "10.76"
$ws.Range("E24").Value = "  -13.16%  "
$ws.Range("E25").Value = "  -15.59%  "
$ws.Range("D26").Value = "5.91"
$ws.Range("E26").Value = "  -3.72%  "
$ws.Range("D27").Value = "3.72"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("E28").Value = "  -11.43%  "
$ws.Range("D29").Value = "11.00"
$ws.Range("E29").Value = "  -12.50%  "
$ws.Range("D30").Value = "8.08"
$ws.Range("E30").Value = "  -12.90%  "
$ws.Range("D31").Value = "634.91"
$ws.Range("E31").Value = "  -10.20%  "
$ws.Range("D32").Value = "27.71"
$ws.Range("E32").Value = "  -12.67%  "
$ws.Range("D33").Value = "6.46"
$ws.Range("E33").Value = "  -15.21%  "
$ws.Range("D34").Value = "10.91"
$ws.Range("E34").Value = "  -10.70%  "
$ws.Range("D35").Value = "57.55"
$ws.Range("E35").Value = "  -9.78%  "
$ws.Range("D38").Value = "34.94"
$ws.Range("E38").Value = "  -17.13%  "
$ws.Range("D39").Value = "0.368"
$ws.Range("E39").Value = "  -11.63%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -9.22%  "
$ws.Range("D42").Value = "2.830.63"
$ws.Range("E42").Value = "  -11.43%  "
$ws.Range("D43").Value = "0.0₃0634"
$ws.Range("E43").Value = "  -16.51%  "
$ws.Range("D44").Value = "2.59"
$ws.Range("E44").Value = "  -20.58%  "
$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  -12.93%  "
$ws.Range("D46").Value = "2.53"
$ws.Range("E46").Value = "  -8.58%  "
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "0.0372"
$ws.Range("E48").Value = "  -10.89%  "
$ws.Range("E49").Value = "  -7.95%  "
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").Value = "  -6.97%  "
$ws.Range("D51").Value = "128.87"
$ws.Range("E51").Value = "  -7.33%  "
